# Re-sort the data rows (everything below the header) alphabetically by
# column A ("file"), ascending. The header row (row 1) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlAscending = 1

$dataRange = $ws.Range("A2:E23")
$sortKey = $ws.Range("A2:A23")
$dataRange.Sort($sortKey, $xlAscending)
